$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 69: record channel-control activity (new shared string "Ovládání kanálu")
$ws.Range("B69").Value = "Ovládání kanálu"

# New rows for recovery after the communication outage
$ws.Range("A73").Value = 42902
$ws.Range("C73").Value = 4

$ws.Range("A74").Value = 42903
$ws.Range("C74").Value = 8

# Move the selection to the newly added row, like the author left it
$ws.Range("A74:C74").Select()
